# Fill the "Definition" column (D) on the "Concepts" sheet with the same
# text currently present in the "Display" column (C), for each concept row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Concepts")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp

for ($r = 2; $r -le $lastRow; $r++) {
    $display = $ws.Cells.Item($r, 3).Value2
    $ws.Cells.Item($r, 4).Value2 = $display
}
